$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C) column date for all existing rows (2-350) from 2023-10-05 to 2023-10-06
$newDate = Get-Date -Year 2023 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2:C350").Value = $newDate

# Row 350 picks up an explicit row height (matches the rest of the sheet)
$ws.Rows.Item(350).RowHeight = 15

# New row 351 data
$ws.Range("A351").Value = "A 47636-2023"

$ws.Range("B351:C351").NumberFormat = "YYYY-MM-DD"
$ws.Range("B351").Value = Get-Date -Year 2023 -Month 10 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("C351").Value = $newDate

$ws.Range("D351").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E351").Value = "NORSJÖ"
$ws.Range("F351").Value = "Holmen skog AB"

$ws.Range("G351").Value = 1
$ws.Range("H351:Q351").Value = 0

$ws.Range("R351").WrapText = $true
